# Recalibration of gold riskDriverType PnL factor calibrations.
# Updates the fitted parameters and abs_epsi_autocorr arrays for the
# AR, SETAR, GARCH, TARCH and AR_TARCH model sheets.

$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.001355084121477351
$ws.Range("B3").Value = 0.784460422979921
$ws.Range("B4").Value = 0.0124336925513482
$ws.Range("B5").Value = "[1.0, 0.24313438848231342, 0.22807726009317775, 0.18786307148781173, 0.17827591184769306, 0.40543785937290117, 0.11673511512741246, 0.1245904039469258, 0.062492425647096096, 0.0650209122358604, 0.017868192808470726, 0.018283071039443387, 0.02997217193708922, 0.03338297220103534, 0.04912243910981346, 0.013995660412526301, 0.019520701918165663, 0.004640502310788266, 0.042267882353372165, 0.041498659335557755]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.03112167765488559
$ws.Range("B5").Value = 0.7400478591889543
$ws.Range("B6").Value = 0.009780676976698474
$ws.Range("B7").Value = 0.02903522103190571
$ws.Range("B8").Value = 0.7308083319954062
$ws.Range("B9").Value = 0.008304740614706463
$ws.Range("B10").Value = "[0.9999999999999998, 0.2237005394391424, 0.24193547414008984, 0.1867812922697022, 0.15051452684525699, 0.2920865368972126, 0.11409170326149035, 0.15433111304588826, 0.09382563982188796, 0.07907618101043881, 0.051600215024133395, 0.0643274123455422, 0.06436865044129719, 0.06953490476527115, 0.08382164172868767, 0.05318614501786818, 0.060788736756055804, 0.035325332113870594, 0.04045628581556891, 0.056161008731325415]"

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.00007755880483153177
$ws.Range("B3").Value = 0.0002786965357638375
$ws.Range("B4").Value = 0.05
$ws.Range("B5").Value = 0.9299999999999999
$ws.Range("B6").Value = "[1.0, 0.03913338072760653, 0.020874618519919307, -0.013772535428039962, -0.022820436896887446, 0.36871159445622675, 0.018098668924043815, 0.035176792099186636, -0.016425063170030743, -0.038720528169120486, -0.06814024934802895, -0.06644768944194926, -0.047730441040245014, -0.013154894591108272, -0.004537972315194989, -0.04249198714239592, -0.03161990509103844, -0.07489218063751679, -0.02810092173319921, -0.025840168682221096]"

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.0008095703757045144
$ws.Range("B3").Value = 0.0003089844554185344
$ws.Range("B4").Value = 0.04987246469385138
$ws.Range("B5").Value = 0.9226407344809797
$ws.Range("B6").Value = "[1.0, 0.03329725706266905, 0.014974607464742384, -0.018467284287051224, -0.027453131000215843, 0.36409805119274957, 0.015468086317117132, 0.03233454061408948, -0.01842866708820073, -0.039908578049522435, -0.06905926144406371, -0.0674786921506437, -0.04865118921439881, -0.013365225752445673, -0.0044437035239168875, -0.04256092918061245, -0.03098647679531979, -0.07477324973963394, -0.028062415154028398, -0.0260257536905864]"
$ws.Range("B7").Value = 0.01507590581183422

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.0008968891137373567
$ws.Range("B3").Value = 0.0004708039690402119
$ws.Range("B4").Value = 0.06739404963471907
$ws.Range("B5").Value = 0.8832691597146155
$ws.Range("B6").Value = "[0.9999999999999999, 0.09071508327853119, 0.07460949812822336, 0.047103457004673145, 0.032956624299207204, 0.26327187823236897, -0.004547668357399099, 0.0036490514608332613, -0.035617818598510695, -0.03408407201077348, -0.07791255090882947, -0.061677421985105105, -0.05312756094979853, -0.02861495981134709, -0.014350636384890772, -0.05434678551296658, -0.03930066488466127, -0.04974780615952934, -0.022230318473961924, -0.017921063674342937]"
$ws.Range("B7").Value = 0.0205642650368734
$ws.Range("B9").Value = 0.7589456554553519
